$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.905.91'
$ws.Range('E2').Value = '  -2.28%  '
$ws.Range('D3').Value = '3.389.46'
$ws.Range('E3').Value = '  -3.32%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.43'
$ws.Range('E5').Value = '  -2.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '126.53'
$ws.Range('E6').Value = '  -5.96%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.390.54'
$ws.Range('E8').Value = '  -3.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.478'
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.31'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E11').Value = '  -4.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.378'
$ws.Range('E12').Value = '  -2.88%  '
$ws.Range('D13').Value = '3.954.87'
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').Value = '3.376.48'
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  -5.45%  '
$ws.Range('D17').Value = '62.872.41'
$ws.Range('E17').Value = '  -2.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.56'
$ws.Range('E18').Value = '  -4.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.32'
$ws.Range('E19').Value = '  -6.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.63'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '372.73'
$ws.Range('E22').Value = '  -5.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.556'
$ws.Range('E23').Value = '  -4.33%  '
$ws.Range('D24').Value = '3.516.34'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.66'
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  -9.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.01'
$ws.Range('E29').Value = '  -5.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.85'
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.12'
$ws.Range('E31').Value = '  -6.75%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  -4.25%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.149'
$ws.Range('E34').Value = '  -5.19%  '
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.411.51'
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.81'
$ws.Range('E36').Value = '  -2.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.43'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '166.05'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.67'
$ws.Range('E39').Value = '  -4.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.49'
$ws.Range('E40').Value = '  -4.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0756'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.79'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.766'
$ws.Range('E44').Value = '  -5.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.24'
$ws.Range('E45').Value = '  -4.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.56'
$ws.Range('E46').Value = '  -6.49%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('E47').Value = '  -6.48%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.50'
$ws.Range('E48').Value = '  -9.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.62'
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('D50').Value = '2.244.33'
$ws.Range('E50').Value = '  -5.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.842'
$ws.Range('E51').Value = '  -7.70%  '
